$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 9).Value = 0.0134664339880133
$ws.Cells.Item(2, 10).Value = 0.0134664339880133
$ws.Cells.Item(2, 13).Value = 5.578493666666667
$ws.Cells.Item(2, 14).Value = 16.735481
$ws.Cells.Item(2, 15).Value = 0.1036332930693284
$ws.Cells.Item(2, 16).Value = 0.1036332930693284
$ws.Cells.Item(2, 17).Value = 0.4654193051036666
$ws.Cells.Item(2, 18).Value = 4.188773745933
$ws.Cells.Item(2, 19).Value = 0.001395570900078547
$ws.Cells.Item(2, 20).Value = 0.001395570900078547
$ws.Cells.Item(3, 9).Value = 0.0134664339880133
$ws.Cells.Item(3, 10).Value = 0.0134664339880133
$ws.Cells.Item(3, 15).Value = 0.06881911773528272
$ws.Cells.Item(3, 16).Value = 0.06881911773528274
$ws.Cells.Item(3, 19).Value = 0.0009267481060955
$ws.Cells.Item(3, 20).Value = 0.0009267481060955002
$ws.Cells.Item(4, 9).Value = 0.0134664339880133
$ws.Cells.Item(4, 10).Value = 0.0134664339880133
$ws.Cells.Item(4, 13).Value = 24.77295966666667
$ws.Cells.Item(4, 14).Value = 74.31887900000001
$ws.Cells.Item(4, 15).Value = 0.4602144490493554
$ws.Cells.Item(4, 16).Value = 0.4602144490493556
$ws.Cells.Item(4, 17).Value = 2.066832797949667
$ws.Cells.Item(4, 18).Value = 18.601495181547
$ws.Cells.Item(4, 19).Value = 0.006197447498453054
$ws.Cells.Item(4, 20).Value = 0.006197447498453058
$ws.Cells.Item(5, 9).Value = 0.0134664339880133
$ws.Cells.Item(5, 10).Value = 0.0134664339880133
$ws.Cells.Item(5, 13).Value = 0.4291063333333334
$ws.Cells.Item(5, 14).Value = 1.287319
$ws.Cells.Item(5, 15).Value = 0.007971632676749163
$ws.Cells.Item(5, 16).Value = 0.007971632676749165
$ws.Cells.Item(5, 17).Value = 0.03580077049633333
$ws.Cells.Item(5, 18).Value = 0.322206934467
$ws.Cells.Item(5, 19).Value = 0.0001073494652181323
$ws.Cells.Item(5, 20).Value = 0.0001073494652181324
$ws.Cells.Item(6, 9).Value = 0.0134664339880133
$ws.Cells.Item(6, 10).Value = 0.0134664339880133
$ws.Cells.Item(6, 13).Value = 19.34413
$ws.Cells.Item(6, 14).Value = 58.03239
$ws.Cells.Item(6, 15).Value = 0.3593615074692841
$ws.Cells.Item(6, 16).Value = 0.3593615074692842
$ws.Cells.Item(6, 17).Value = 1.61390011003
$ws.Cells.Item(6, 18).Value = 14.52510099027
$ws.Cells.Item(6, 19).Value = 0.004839318018168062
$ws.Cells.Item(6, 20).Value = 0.004839318018168064
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 6).Value = 1
$ws.Cells.Item(7, 7).Value = 4.699187666666667
$ws.Cells.Item(7, 8).Value = 14.097563
$ws.Cells.Item(7, 9).Value = 0.758486659760196
$ws.Cells.Item(7, 10).Value = 0.758486659760196
$ws.Cells.Item(7, 13).Value = 5.578493666666667
$ws.Cells.Item(7, 14).Value = 16.735481
$ws.Cells.Item(7, 15).Value = 0.1036332930693284
$ws.Cells.Item(7, 16).Value = 0.1036332930693284
$ws.Cells.Item(7, 17).Value = 26.21438863697811
$ws.Cells.Item(7, 18).Value = 235.929497732803
$ws.Cells.Item(7, 19).Value = 0.07860447030010437
$ws.Cells.Item(7, 20).Value = 0.07860447030010438
$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 6).Value = 1
$ws.Cells.Item(8, 7).Value = 4.699187666666667
$ws.Cells.Item(8, 8).Value = 14.097563
$ws.Cells.Item(8, 9).Value = 0.758486659760196
$ws.Cells.Item(8, 10).Value = 0.758486659760196
$ws.Cells.Item(8, 15).Value = 0.06881911773528272
$ws.Cells.Item(8, 16).Value = 0.06881911773528274
$ws.Cells.Item(8, 17).Value = 17.40802636426678
$ws.Cells.Item(8, 18).Value = 156.672237278401
$ws.Cells.Item(8, 19).Value = 0.05219838273867826
$ws.Cells.Item(8, 20).Value = 0.05219838273867827
$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 6).Value = 1
$ws.Cells.Item(9, 7).Value = 4.699187666666667
$ws.Cells.Item(9, 8).Value = 14.097563
$ws.Cells.Item(9, 9).Value = 0.758486659760196
$ws.Cells.Item(9, 10).Value = 0.758486659760196
$ws.Cells.Item(9, 13).Value = 24.77295966666667
$ws.Cells.Item(9, 14).Value = 74.31887900000001
$ws.Cells.Item(9, 15).Value = 0.4602144490493554
$ws.Cells.Item(9, 16).Value = 0.4602144490493556
$ws.Cells.Item(9, 17).Value = 116.4127865324308
$ws.Cells.Item(9, 18).Value = 1047.715078791877
$ws.Cells.Item(9, 19).Value = 0.3490665202328245
$ws.Cells.Item(9, 20).Value = 0.3490665202328246
$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 6).Value = 1
$ws.Cells.Item(10, 7).Value = 4.699187666666667
$ws.Cells.Item(10, 8).Value = 14.097563
$ws.Cells.Item(10, 9).Value = 0.758486659760196
$ws.Cells.Item(10, 10).Value = 0.758486659760196
$ws.Cells.Item(10, 13).Value = 0.4291063333333334
$ws.Cells.Item(10, 14).Value = 1.287319
$ws.Cells.Item(10, 15).Value = 0.007971632676749163
$ws.Cells.Item(10, 16).Value = 0.007971632676749165
$ws.Cells.Item(10, 17).Value = 2.016451189288556
$ws.Cells.Item(10, 18).Value = 18.148060703597
$ws.Cells.Item(10, 19).Value = 0.006046377041822703
$ws.Cells.Item(10, 20).Value = 0.006046377041822705
$ws.Cells.Item(11, 5).Value = 3
$ws.Cells.Item(11, 6).Value = 1
$ws.Cells.Item(11, 7).Value = 4.699187666666667
$ws.Cells.Item(11, 8).Value = 14.097563
$ws.Cells.Item(11, 9).Value = 0.758486659760196
$ws.Cells.Item(11, 10).Value = 0.758486659760196
$ws.Cells.Item(11, 13).Value = 19.34413
$ws.Cells.Item(11, 14).Value = 58.03239
$ws.Cells.Item(11, 15).Value = 0.3593615074692841
$ws.Cells.Item(11, 16).Value = 0.3593615074692842
$ws.Cells.Item(11, 17).Value = 90.90169711839667
$ws.Cells.Item(11, 18).Value = 818.11527406557
$ws.Cells.Item(11, 19).Value = 0.272570909446766
$ws.Cells.Item(11, 20).Value = 0.2725709094467661
$ws.Cells.Item(12, 7).Value = 0.1492686666666667
$ws.Cells.Item(12, 8).Value = 0.447806
$ws.Cells.Item(12, 9).Value = 0.02409316256721636
$ws.Cells.Item(12, 10).Value = 0.02409316256721636
$ws.Cells.Item(12, 13).Value = 5.578493666666667
$ws.Cells.Item(12, 14).Value = 16.735481
$ws.Cells.Item(12, 15).Value = 0.1036332930693284
$ws.Cells.Item(12, 16).Value = 0.1036332930693284
$ws.Cells.Item(12, 17).Value = 0.8326943116317779
$ws.Cells.Item(12, 18).Value = 7.494248804686
$ws.Cells.Item(12, 19).Value = 0.002496853777295306
$ws.Cells.Item(12, 20).Value = 0.002496853777295306
$ws.Cells.Item(13, 7).Value = 0.1492686666666667
$ws.Cells.Item(13, 8).Value = 0.447806
$ws.Cells.Item(13, 9).Value = 0.02409316256721636
$ws.Cells.Item(13, 10).Value = 0.02409316256721636
$ws.Cells.Item(13, 15).Value = 0.06881911773528272
$ws.Cells.Item(13, 16).Value = 0.06881911773528274
$ws.Cells.Item(13, 17).Value = 0.5529621434624445
$ws.Cells.Item(13, 18).Value = 4.976659291162
$ws.Cells.Item(13, 19).Value = 0.001658070191328569
$ws.Cells.Item(13, 20).Value = 0.00165807019132857
$ws.Cells.Item(14, 7).Value = 0.1492686666666667
$ws.Cells.Item(14, 8).Value = 0.447806
$ws.Cells.Item(14, 9).Value = 0.02409316256721636
$ws.Cells.Item(14, 10).Value = 0.02409316256721636
$ws.Cells.Item(14, 13).Value = 24.77295966666667
$ws.Cells.Item(14, 14).Value = 74.31887900000001
$ws.Cells.Item(14, 15).Value = 0.4602144490493554
$ws.Cells.Item(14, 16).Value = 0.4602144490493556
$ws.Cells.Item(14, 17).Value = 3.697826658830445
$ws.Cells.Item(14, 18).Value = 33.28043992947401
$ws.Cells.Item(14, 19).Value = 0.01108802153672803
$ws.Cells.Item(14, 20).Value = 0.01108802153672804
$ws.Cells.Item(15, 7).Value = 0.1492686666666667
$ws.Cells.Item(15, 8).Value = 0.447806
$ws.Cells.Item(15, 9).Value = 0.02409316256721636
$ws.Cells.Item(15, 10).Value = 0.02409316256721636
$ws.Cells.Item(15, 13).Value = 0.4291063333333334
$ws.Cells.Item(15, 14).Value = 1.287319
$ws.Cells.Item(15, 15).Value = 0.007971632676749163
$ws.Cells.Item(15, 16).Value = 0.007971632676749165
$ws.Cells.Item(15, 17).Value = 0.0640521302348889
$ws.Cells.Item(15, 18).Value = 0.5764691721140001
$ws.Cells.Item(15, 19).Value = 0.0001920618420070517
$ws.Cells.Item(15, 20).Value = 0.0001920618420070517
$ws.Cells.Item(16, 7).Value = 0.1492686666666667
$ws.Cells.Item(16, 8).Value = 0.447806
$ws.Cells.Item(16, 9).Value = 0.02409316256721636
$ws.Cells.Item(16, 10).Value = 0.02409316256721636
$ws.Cells.Item(16, 13).Value = 19.34413
$ws.Cells.Item(16, 14).Value = 58.03239
$ws.Cells.Item(16, 15).Value = 0.3593615074692841
$ws.Cells.Item(16, 16).Value = 0.3593615074692842
$ws.Cells.Item(16, 17).Value = 2.887472492926667
$ws.Cells.Item(16, 18).Value = 25.98725243634
$ws.Cells.Item(16, 19).Value = 0.0086581552198574
$ws.Cells.Item(16, 20).Value = 0.0086581552198574
$ws.Cells.Item(17, 5).Value = 3
$ws.Cells.Item(17, 6).Value = 1
$ws.Cells.Item(17, 7).Value = 1.263591
$ws.Cells.Item(17, 8).Value = 3.790773
$ws.Cells.Item(17, 9).Value = 0.2039537436845743
$ws.Cells.Item(17, 10).Value = 0.2039537436845743
$ws.Cells.Item(17, 13).Value = 5.578493666666667
$ws.Cells.Item(17, 14).Value = 16.735481
$ws.Cells.Item(17, 15).Value = 0.1036332930693284
$ws.Cells.Item(17, 16).Value = 0.1036332930693284
$ws.Cells.Item(17, 17).Value = 7.048934390757001
$ws.Cells.Item(17, 18).Value = 63.44040951681301
$ws.Cells.Item(17, 19).Value = 0.02113639809185017
$ws.Cells.Item(17, 20).Value = 0.02113639809185017
$ws.Cells.Item(18, 5).Value = 3
$ws.Cells.Item(18, 6).Value = 1
$ws.Cells.Item(18, 7).Value = 1.263591
$ws.Cells.Item(18, 8).Value = 3.790773
$ws.Cells.Item(18, 9).Value = 0.2039537436845743
$ws.Cells.Item(18, 10).Value = 0.2039537436845743
$ws.Cells.Item(18, 15).Value = 0.06881911773528272
$ws.Cells.Item(18, 16).Value = 0.06881911773528274
$ws.Cells.Item(18, 17).Value = 4.680942112119
$ws.Cells.Item(18, 18).Value = 42.128479009071
$ws.Cells.Item(18, 19).Value = 0.01403591669918039
$ws.Cells.Item(18, 20).Value = 0.01403591669918039
$ws.Cells.Item(19, 5).Value = 3
$ws.Cells.Item(19, 6).Value = 1
$ws.Cells.Item(19, 7).Value = 1.263591
$ws.Cells.Item(19, 8).Value = 3.790773
$ws.Cells.Item(19, 9).Value = 0.2039537436845743
$ws.Cells.Item(19, 10).Value = 0.2039537436845743
$ws.Cells.Item(19, 13).Value = 24.77295966666667
$ws.Cells.Item(19, 14).Value = 74.31887900000001
$ws.Cells.Item(19, 15).Value = 0.4602144490493554
$ws.Cells.Item(19, 16).Value = 0.4602144490493556
$ws.Cells.Item(19, 17).Value = 31.30288887816301
$ws.Cells.Item(19, 18).Value = 281.725999903467
$ws.Cells.Item(19, 19).Value = 0.0938624597813498
$ws.Cells.Item(19, 20).Value = 0.09386245978134983
$ws.Cells.Item(20, 5).Value = 3
$ws.Cells.Item(20, 6).Value = 1
$ws.Cells.Item(20, 7).Value = 1.263591
$ws.Cells.Item(20, 8).Value = 3.790773
$ws.Cells.Item(20, 9).Value = 0.2039537436845743
$ws.Cells.Item(20, 10).Value = 0.2039537436845743
$ws.Cells.Item(20, 13).Value = 0.4291063333333334
$ws.Cells.Item(20, 14).Value = 1.287319
$ws.Cells.Item(20, 15).Value = 0.007971632676749163
$ws.Cells.Item(20, 16).Value = 0.007971632676749165
$ws.Cells.Item(20, 17).Value = 0.5422149008430001
$ws.Cells.Item(20, 18).Value = 4.879934107587
$ws.Cells.Item(20, 19).Value = 0.001625844327701275
$ws.Cells.Item(20, 20).Value = 0.001625844327701276
$ws.Cells.Item(21, 5).Value = 3
$ws.Cells.Item(21, 6).Value = 1
$ws.Cells.Item(21, 7).Value = 1.263591
$ws.Cells.Item(21, 8).Value = 3.790773
$ws.Cells.Item(21, 9).Value = 0.2039537436845743
$ws.Cells.Item(21, 10).Value = 0.2039537436845743
$ws.Cells.Item(21, 13).Value = 19.34413
$ws.Cells.Item(21, 14).Value = 58.03239
$ws.Cells.Item(21, 15).Value = 0.3593615074692841
$ws.Cells.Item(21, 16).Value = 0.3593615074692842
$ws.Cells.Item(21, 17).Value = 24.44306857083
$ws.Cells.Item(21, 18).Value = 219.98761713747
$ws.Cells.Item(21, 19).Value = 0.0732931247844926
$ws.Cells.Item(21, 20).Value = 0.0732931247844926
